$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows with updated calibration costs, each applied uniformly across columns J:AS (cols 10-45)
$updates = @{
    100 = 14447.15351
    101 = 20769.84987
    102 = 769169.9638
    103 = 5240.939598
    104 = 11113.19501
    105 = 1161.553908
    106 = 84608.69602
    107 = 27054.76164
    114 = 20.64299912
    115 = 6633973.386
}

foreach ($r in $updates.Keys) {
    $value = $updates[$r]
    $rng = $ws.Range($ws.Cells.Item($r, 10), $ws.Cells.Item($r, 45))
    $rng.Value = $value
}
